$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 308.31146
$ws.Range("J17").Value = 308.31146
$ws.Range("L17").Value = 924.93438
$ws.Range("N17").Value = -1260.93438
$ws.Range("H64").Value = 128112.375
$ws.Range("I64").Value = 335466.66
$ws.Range("J64").Value = 3699.8
$ws.Range("K64").Value = 335466.66
$ws.Range("L64").Value = 3699.8
$ws.Range("M64").Value = -335218.66
$ws.Range("N64").Value = -4195.8
$ws.Range("H67").Value = 128112.375
$ws.Range("I67").Value = 335466.66
$ws.Range("J67").Value = 3699.8
$ws.Range("K67").Value = 335466.66
$ws.Range("L67").Value = 3699.8
$ws.Range("M67").Value = -334608.66
$ws.Range("N67").Value = -5415.8
$ws.Range("H92").Value = 1097.55
$ws.Range("I92").Value = 1142.8889
$ws.Range("J92").Value = 689.5
$ws.Range("K92").Value = 1142.8889
$ws.Range("L92").Value = 689.5
$ws.Range("M92").Value = 105.1111000000001
$ws.Range("N92").Value = -3185.5
$ws.Range("H137").Value = 1515.3334
$ws.Range("I137").Value = 1550.6666
$ws.Range("J137").Value = 1497.6666
$ws.Range("K137").Value = 4651.9998
$ws.Range("L137").Value = 4492.9998
$ws.Range("M137").Value = -2101.9998
$ws.Range("N137").Value = -9592.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 37726.24
$ws.Range("I32").Value = 6082.231
$ws.Range("K32").Value = 6082.231
$ws.Range("M32").Value = -5795.231
$ws.Range("H74").Value = 2101.7112
$ws.Range("I74").Value = 1588.9259
$ws.Range("J74").Value = 2870.889
$ws.Range("K74").Value = 1588.9259
$ws.Range("L74").Value = 2870.889
$ws.Range("M74").Value = -714.9259
$ws.Range("N74").Value = -4618.889
$ws.Range("H77").Value = 2101.7112
$ws.Range("I77").Value = 1588.9259
$ws.Range("J77").Value = 2870.889
$ws.Range("K77").Value = 7944.6295
$ws.Range("L77").Value = 14354.445
$ws.Range("M77").Value = -3576.6295
$ws.Range("N77").Value = -23090.445
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 999.8333
$ws.Range("I16").Value = 799.6667
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 799.6667
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -512.6667
$ws.Range("N16").Value = -1774
$ws.Range("H31").Value = 26591.195
$ws.Range("I31").Value = 41396.12
$ws.Range("J31").Value = 3458.5
$ws.Range("K31").Value = 41396.12
$ws.Range("L31").Value = 3458.5
$ws.Range("M31").Value = -41101.12
$ws.Range("N31").Value = -4048.5
$ws.Range("H34").Value = 26591.195
$ws.Range("I34").Value = 41396.12
$ws.Range("J34").Value = 3458.5
$ws.Range("K34").Value = 41396.12
$ws.Range("L34").Value = 3458.5
$ws.Range("M34").Value = -41194.12
$ws.Range("N34").Value = -3862.5
$ws.Range("H50").Value = 14428
$ws.Range("J50").Value = 14428
$ws.Range("L50").Value = 14428
$ws.Range("N50").Value = -15678
$ws.Range("H51").Value = 7917
$ws.Range("J51").Value = 7897.778
$ws.Range("L51").Value = 7897.778
$ws.Range("N51").Value = -9369.778
$ws.Range("H58").Value = 10537.036
$ws.Range("I58").Value = 1658.7858
$ws.Range("K58").Value = 1658.7858
$ws.Range("M58").Value = -1455.7858
$ws.Range("H61").Value = 7917
$ws.Range("J61").Value = 7897.778
$ws.Range("L61").Value = 7897.778
$ws.Range("N61").Value = -8593.778
$ws.Range("H68").Value = 14356.357
$ws.Range("J68").Value = 14356.357
$ws.Range("L68").Value = 14356.357
$ws.Range("N68").Value = -15854.357
$ws.Range("H71").Value = 14356.357
$ws.Range("J71").Value = 14356.357
$ws.Range("L71").Value = 43069.071
$ws.Range("N71").Value = -50557.071
$ws.Range("H74").Value = 23361.637
$ws.Range("J74").Value = 23361.637
$ws.Range("L74").Value = 23361.637
$ws.Range("N74").Value = -25109.637
$ws.Range("H77").Value = 23361.637
$ws.Range("J77").Value = 23361.637
$ws.Range("L77").Value = 70084.91099999999
$ws.Range("N77").Value = -78820.91099999999
$ws.Range("H99").Value = 8746.4375
$ws.Range("I99").Value = 2734
$ws.Range("J99").Value = 11479.363
$ws.Range("K99").Value = 2734
$ws.Range("L99").Value = 11479.363
$ws.Range("M99").Value = -1236
$ws.Range("N99").Value = -14475.363
$ws.Range("H113").Value = 999.8333
$ws.Range("I113").Value = 799.6667
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 799.6667
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1370.3333
$ws.Range("N113").Value = -5540
$ws.Range("H126").Value = 8746.4375
$ws.Range("I126").Value = 2734
$ws.Range("J126").Value = 11479.363
$ws.Range("K126").Value = 8202
$ws.Range("L126").Value = 34438.089
$ws.Range("M126").Value = -5732
$ws.Range("N126").Value = -39378.089
$ws.Range("H136").Value = 10537.036
$ws.Range("I136").Value = 1658.7858
$ws.Range("K136").Value = 4976.357400000001
$ws.Range("M136").Value = -2426.357400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 2050
$ws.Range("J41").Value = 3600
$ws.Range("L41").Value = 10800
$ws.Range("N41").Value = -11476
$ws.Range("H103").Value = 1358.2222
$ws.Range("I103").Value = 912.3333
$ws.Range("J103").Value = 2250
$ws.Range("K103").Value = 2736.9999
$ws.Range("L103").Value = 6750
$ws.Range("M103").Value = -1857.9999
$ws.Range("N103").Value = -8508
$ws.Range("H121").Value = 2622.7778
$ws.Range("I121").Value = 364.875
$ws.Range("J121").Value = 4429.1
$ws.Range("K121").Value = 1094.625
$ws.Range("L121").Value = 13287.3
$ws.Range("M121").Value = 215.375
$ws.Range("N121").Value = -15907.3
$ws.Range("H122").Value = 7314.467
$ws.Range("J122").Value = 17649.666
$ws.Range("L122").Value = 158846.994
$ws.Range("N122").Value = -163746.994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2081.2144
$ws.Range("J113").Value = 1631.7778
$ws.Range("L113").Value = 1631.7778
$ws.Range("N113").Value = -5971.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6750
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 6750
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 6750
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -7126
$ws.Range("H100").Value = 2121
$ws.Range("I100").Value = 2033.3334
$ws.Range("J100").Value = 2186.75
$ws.Range("K100").Value = 2033.3334
$ws.Range("L100").Value = 2186.75
$ws.Range("M100").Value = -1492.3334
$ws.Range("N100").Value = -3268.75
